$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    # Force the literal text into the cell without Excel re-parsing it as a
    # number/date, then drop the temporary text format so the cell keeps its
    # original (default) style, matching the source diff exactly.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

Set-TextValue $ws.Range("D2") "53.637.56"
Set-TextValue $ws.Range("E2") "  -4.64%  "
Set-TextValue $ws.Range("D3") "2.225.48"
Set-TextValue $ws.Range("E3") "  -6.15%  "
Set-TextValue $ws.Range("D4") "1.00"
Set-TextValue $ws.Range("E4") "  -0.11%  "
Set-TextValue $ws.Range("D5") "482.75"
Set-TextValue $ws.Range("E5") "  -3.68%  "
Set-TextValue $ws.Range("D6") "125.10"
Set-TextValue $ws.Range("E6") "  -3.48%  "
Set-TextValue $ws.Range("D7") "1.00"
Set-TextValue $ws.Range("E7") "  +0.00%  "
Set-TextValue $ws.Range("D8") "0.516"
Set-TextValue $ws.Range("E8") "  -5.37%  "
Set-TextValue $ws.Range("D9") "2.237.02"
Set-TextValue $ws.Range("E9") "  -5.74%  "
Set-TextValue $ws.Range("D10") "0.0914"
Set-TextValue $ws.Range("E10") "  -6.83%  "
Set-TextValue $ws.Range("E11") "  -1.51%  "
Set-TextValue $ws.Range("D12") "4.67"
Set-TextValue $ws.Range("E12") "  -2.21%  "
Set-TextValue $ws.Range("D13") "0.314"
Set-TextValue $ws.Range("E13") "  -3.13%  "
Set-TextValue $ws.Range("D14") "2.623.90"
Set-TextValue $ws.Range("E14") "  -6.21%  "
Set-TextValue $ws.Range("D15") "20.97"
Set-TextValue $ws.Range("E15") "  -2.39%  "
Set-TextValue $ws.Range("D16") "53.568.37"
Set-TextValue $ws.Range("E16") "  -4.76%  "
Set-TextValue $ws.Range("E17") "  -3.76%  "
Set-TextValue $ws.Range("D18") "2.220.62"
Set-TextValue $ws.Range("E18") "  -6.20%  "
Set-TextValue $ws.Range("D19") "9.54"
Set-TextValue $ws.Range("E19") "  -4.71%  "
Set-TextValue $ws.Range("D20") "3.95"
Set-TextValue $ws.Range("E20") "  -2.12%  "
Set-TextValue $ws.Range("D21") "297.95"
Set-TextValue $ws.Range("E21") "  -2.88%  "
Set-TextValue $ws.Range("E22") "  -2.73%  "
Set-TextValue $ws.Range("E23") "  +0.32%  "
Set-TextValue $ws.Range("E24") "  -3.27%  "
Set-TextValue $ws.Range("E25") "  +0.19%  "
Set-TextValue $ws.Range("E26") "  -1.49%  "
Set-TextValue $ws.Range("E27") "  -3.64%  "
Set-TextValue $ws.Range("E28") "  -3.75%  "
Set-TextValue $ws.Range("D29") "170.11"
Set-TextValue $ws.Range("E29") "  -0.82%  "
Set-TextValue $ws.Range("D30") "0.0₃0675"
Set-TextValue $ws.Range("E30") "  -5.51%  "
Set-TextValue $ws.Range("E31") "  -4.12%  "
Set-TextValue $ws.Range("E32") "  -0.12%  "
Set-TextValue $ws.Range("D33") "0.998"
Set-TextValue $ws.Range("E33") "  -0.03%  "
Set-TextValue $ws.Range("D34") "5.69"
Set-TextValue $ws.Range("E34") "  -0.94%  "
Set-TextValue $ws.Range("E35") "  -3.70%  "
Set-TextValue $ws.Range("D36") "17.37"
Set-TextValue $ws.Range("E36") "  -1.01%  "
Set-TextValue $ws.Range("E37") "  -2.54%  "
Set-TextValue $ws.Range("D38") "0.829"
Set-TextValue $ws.Range("E38") "  +4.37%  "
Set-TextValue $ws.Range("E39") "  -5.23%  "
Set-TextValue $ws.Range("D40") "35.86"
Set-TextValue $ws.Range("E40") "  -0.57%  "
Set-TextValue $ws.Range("E41") "  -1.02%  "
Set-TextValue $ws.Range("E42") "  -1.85%  "
Set-TextValue $ws.Range("E43") "  -2.51%  "
Set-TextValue $ws.Range("D44") "122.54"
Set-TextValue $ws.Range("E44") "  -6.08%  "
Set-TextValue $ws.Range("D45") "4.65"
Set-TextValue $ws.Range("E45") "  -1.47%  "
Set-TextValue $ws.Range("E46") "  -3.34%  "
Set-TextValue $ws.Range("D47") "0.534"
Set-TextValue $ws.Range("E47") "  -4.84%  "
Set-TextValue $ws.Range("E48") "  -2.46%  "
Set-TextValue $ws.Range("D49") "228.36"
Set-TextValue $ws.Range("E49") "  -4.83%  "
Set-TextValue $ws.Range("D50") "0.0200"
Set-TextValue $ws.Range("E50") "  -3.27%  "
Set-TextValue $ws.Range("D51") "15.94"
Set-TextValue $ws.Range("E51") "  -5.61%  "
